# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh).
# Forces each written value to remain plain text (matching the workbook's
# existing inline-string cells) even when the new value happens to look
# like a number (e.g. "1.00"), and resets the cell style afterwards so no
# stray number-format style is left behind.
function Set-CellText {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "63.323.00"
Set-CellText $ws "E2" "  -2.65%  "
Set-CellText $ws "D3" "3.120.45"
Set-CellText $ws "E3" "  -0.98%  "
Set-CellText $ws "E4" "  -0.01%  "
Set-CellText $ws "D5" "557.71"
Set-CellText $ws "E5" "  -1.07%  "
Set-CellText $ws "D6" "139.47"
Set-CellText $ws "E6" "  -6.59%  "
Set-CellText $ws "D7" "1.00"
Set-CellText $ws "E7" "  +0.12%  "
Set-CellText $ws "D8" "3.113.93"
Set-CellText $ws "D9" "0.498"
Set-CellText $ws "E9" "  -0.73%  "
Set-CellText $ws "B10" "Dogecoin"
Set-CellText $ws "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-CellText $ws "D10" "0.161"
Set-CellText $ws "E10" "  -0.52%  "
Set-CellText $ws "B11" "Toncoin"
Set-CellText $ws "C11" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D11" "6.65"
Set-CellText $ws "E11" "  -1.18%  "
Set-CellText $ws "E12" "  -1.45%  "
Set-CellText $ws "D13" "35.43"
Set-CellText $ws "E13" "  -5.63%  "
Set-CellText $ws "E14" "  -2.87%  "
Set-CellText $ws "D15" "3.633.48"
Set-CellText $ws "E15" "  -0.80%  "
Set-CellText $ws "D16" "63.373.12"
Set-CellText $ws "E16" "  -2.56%  "
Set-CellText $ws "E17" "  -0.62%  "
Set-CellText $ws "D18" "3.122.08"
Set-CellText $ws "E18" "  -1.04%  "
Set-CellText $ws "D19" "509.62"
Set-CellText $ws "E19" "  -0.69%  "
Set-CellText $ws "D20" "6.75"
Set-CellText $ws "E20" "  -1.35%  "
Set-CellText $ws "D21" "13.66"
Set-CellText $ws "E21" "  -3.08%  "
Set-CellText $ws "D22" "0.713"
Set-CellText $ws "E22" "  +0.75%  "
Set-CellText $ws "E23" "  -1.73%  "
Set-CellText $ws "D24" "12.50"
Set-CellText $ws "E24" "  -3.03%  "
Set-CellText $ws "D25" "78.39"
Set-CellText $ws "E25" "  -1.37%  "
Set-CellText $ws "E26" "  +0.06%  "
Set-CellText $ws "D27" "2.78"
Set-CellText $ws "E27" "  -0.94%  "
Set-CellText $ws "D28" "8.35"
Set-CellText $ws "E28" "  -3.44%  "
Set-CellText $ws "D29" "1.00"
Set-CellText $ws "E29" "  -0.12%  "
Set-CellText $ws "E30" "  -9.12%  "
Set-CellText $ws "D31" "26.49"
Set-CellText $ws "E31" "  -1.04%  "
Set-CellText $ws "E32" "  -7.18%  "
Set-CellText $ws "E33" "  -2.13%  "
Set-CellText $ws "D34" "59.48"
Set-CellText $ws "E34" "  +11.85%  "
Set-CellText $ws "D35" "533.54"
Set-CellText $ws "E35" "  -11.48%  "
Set-CellText $ws "E36" "  -1.72%  "
Set-CellText $ws "E37" "  -6.83%  "
Set-CellText $ws "D38" "0.0418"
Set-CellText $ws "E38" "  -3.28%  "
Set-CellText $ws "D39" "0.0801"
Set-CellText $ws "E39" "  -2.65%  "
Set-CellText $ws "D40" "3.080.72"
Set-CellText $ws "E40" "  +0.64%  "
Set-CellText $ws "E41" "  -2.71%  "
Set-CellText $ws "D42" "2.75"
Set-CellText $ws "E42" "  -7.86%  "
Set-CellText $ws "D43" "8.15"
Set-CellText $ws "E43" "  -2.45%  "
Set-CellText $ws "D44" "0.257"
Set-CellText $ws "E44" "  -0.57%  "
Set-CellText $ws "D46" "2.08"
Set-CellText $ws "E46" "  -5.47%  "
Set-CellText $ws "D47" "122.68"
Set-CellText $ws "E47" "  +0.58%  "
Set-CellText $ws "D48" "24.47"
Set-CellText $ws "E48" "  -6.27%  "
Set-CellText $ws "E49" "  -2.19%  "
Set-CellText $ws "B50" "CoreDAO"
Set-CellText $ws "C50" "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-CellText $ws "D50" "2.46"
Set-CellText $ws "E50" "  +62.72%  "
Set-CellText $ws "B51" "PEPE"
Set-CellText $ws "C51" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText $ws "D51" "0.0₃0513"
Set-CellText $ws "E51" "  -6.44%  "
